$d = $word.ActiveDocument

$ids = @("p049v_1", "p049v_2", "p049v_3")

foreach ($id in $ids) {
    $target = "<id>" + $id + "</id>"
    $rng = $d.Content
    $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $target, 2)
}
